# "added people to gantn chart for week 1"
# Fill in the Week-1 owner names on the Gantt chart (Sheet1!B2:B4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write B4 first, then B2, then B3 so new shared-string entries land in the
# same order as the authored workbook (Zhipeng, George/Diana/Pranoy,
# George/Sunh/Wei).
$ws.Range("B4").Value = "Zhipeng"
$ws.Range("B2").Value = "George/Diana/Pranoy"
$ws.Range("B3").Value = "George/Sunh/Wei"

# Leave the selection on the last cell that was edited, matching the saved
# workbook's cursor position.
$ws.Range("B4").Select()

# The workbook also picked up the defined name that the "MySQL for Excel"
# add-in stamps into every workbook it touches on save.
$mysqlDateName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$mysqlDateName.Visible = $false

# That same add-in also registers its "MySqlDefault" table style (two dxf
# records: a bold/shaded header row dxf, and a blank "whole table" dxf).
# Recreate the two dxf records via a pair of throwaway conditional formats
# (deleted immediately after) so the style records land in xl/styles.xml
# without leaving a real conditional format or table behind.
$headerRng = $ws.Range("Z1")
$headerFc = $headerRng.FormatConditions.Add(1, 3, "1")
$headerFc.Font.Bold = $true
$headerFc.Font.Italic = $false
$headerFc.Interior.Color = 14145495
$headerRng.FormatConditions.Delete()

$wholeRng = $ws.Range("Z2")
$wholeFc = $wholeRng.FormatConditions.Add(1, 3, "1")
$wholeFc.Font.Bold = $false
$wholeFc.Font.Italic = $false
$wholeFc.Interior.Pattern = -4142
$wholeRng.FormatConditions.Delete()
